$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Agregar info en experiencia y educacion: fecha de inicio en B2 con formato de fecha
$ws.Range("B2").HorizontalAlignment = 1
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = 44582

# Mover la seleccion activa a A3
$ws.Range("A3").Select()
